# add genders in batumi
# Adds a new "2023" year column (T) to the hotels & restaurants economic
# indicators sheet, mirroring the existing per-year columns B..S, and widens
# the data columns to fit the now-larger B:S range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column T: copy each row's formatting from column S, then set the
#     2023 figure (mirrors how the 2022 column itself looks) ---
$pairs = @(
    @{ Src = "S3";  Dst = "T3";  Val = 2023 },
    @{ Src = "S4";  Dst = "T4";  Val = 455.5 },
    @{ Src = "S5";  Dst = "T5";  Val = 483.8 },
    @{ Src = "S6";  Dst = "T6";  Val = 5833 },
    @{ Src = "S7";  Dst = "T7";  Val = 5749 },
    @{ Src = "S8";  Dst = "T8";  Val = 1408.7 },
    @{ Src = "S9";  Dst = "T9";  Val = 253.6 },
    @{ Src = "S10"; Dst = "T10"; Val = 99.1 },
    @{ Src = "S11"; Dst = "T11"; Val = 230.1 },
    @{ Src = "S12"; Dst = "T12"; Val = 233.6 },
    @{ Src = "S13"; Dst = "T13"; Val = 260.3 },
    @{ Src = "S14"; Dst = "T14"; Val = 0.1 }
)

foreach ($p in $pairs) {
    $ws.Range($p.Src).Copy()
    $ws.Range($p.Dst).PasteSpecial(-4122)  # xlPasteFormats
    $ws.Range($p.Dst).Value = $p.Val
}

# --- Column widths: B:S (now 18 columns, was B:P/16) get a uniform width ---
$ws.Range("B1:S1").EntireColumn.ColumnWidth = 8

# --- Restore the authored selection state ---
$ws.Range("W10").Select()
